# Set up sustained reentry circuit
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Path")

# Update path conduction/refractory values
$ws.Range("E2").Value = 30
$ws.Range("J2").Value = 30

$ws.Range("E3").Value = 30
$ws.Range("J3").Value = 30

$ws.Range("E4").Value = 1200
$ws.Range("J4").Value = 1200

$ws.Range("E5").Value = 1200
$ws.Range("J5").Value = 1200

# Update the active selection on the Path sheet
$ws.Activate()
$ws.Range("L14").Select()
